$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Penk"
$ws.Range("C2").Value = "Mrgprb2"
$ws.Range("D2").Value = "Neutrophils"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.049583
$ws.Range("H2").Value = 0.148749
$ws.Range("I2").Value = 0.003300257029073341
$ws.Range("J2").Value = 0.003300257029073341
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.4500266666666666
$ws.Range("N2").Value = 1.35008
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.02231367221333333
$ws.Range("R2").Value = 0.20082304992
$ws.Range("S2").Value = 0.003300257029073341
$ws.Range("T2").Value = 0.003300257029073341

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Penk"
$ws.Range("C3").Value = "Mrgprb2"
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 13.14903766666666
$ws.Range("H3").Value = 39.44711299999999
$ws.Range("I3").Value = 0.8752032750129437
$ws.Range("J3").Value = 0.8752032750129437
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.4500266666666666
$ws.Range("N3").Value = 1.35008
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 5.917417591004443
$ws.Range("R3").Value = 53.25675831903999
$ws.Range("S3").Value = 0.8752032750129437
$ws.Range("T3").Value = 0.8752032750129437

$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Penk"
$ws.Range("C4").Value = "Mrgprb2"
$ws.Range("D4").Value = "Neutrophils"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.411192
$ws.Range("H4").Value = 1.233576
$ws.Range("I4").Value = 0.02736904358951103
$ws.Range("J4").Value = 0.02736904358951103
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.4500266666666666
$ws.Range("N4").Value = 1.35008
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.18504736512
$ws.Range("R4").Value = 1.66542628608
$ws.Range("S4").Value = 0.02736904358951103
$ws.Range("T4").Value = 0.02736904358951103

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Penk"
$ws.Range("C5").Value = "Mrgprb2"
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.860481
$ws.Range("H5").Value = 2.581443
$ws.Range("I5").Value = 0.05727383314107775
$ws.Range("J5").Value = 0.05727383314107775
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.4500266666666666
$ws.Range("N5").Value = 1.35008
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.38723939616
$ws.Range("R5").Value = 3.48515456544
$ws.Range("S5").Value = 0.05727383314107775
$ws.Range("T5").Value = 0.05727383314107775

$ws.Range("A6").Value = "Neutrophils"
$ws.Range("B6").Value = "Penk"
$ws.Range("C6").Value = "Mrgprb2"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2591163333333333
$ws.Range("H6").Value = 0.777349
$ws.Range("I6").Value = 0.0172468487270041
$ws.Range("J6").Value = 0.0172468487270041
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.4500266666666666
$ws.Range("N6").Value = 1.35008
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 0.1166092597688889
$ws.Range("R6").Value = 1.04948333792
$ws.Range("S6").Value = 0.0172468487270041
$ws.Range("T6").Value = 0.0172468487270041

$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Penk"
$ws.Range("C7").Value = "Mrgprb2"
$ws.Range("D7").Value = "Neutrophils"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2945713333333334
$ws.Range("H7").Value = 0.883714
$ws.Range("I7").Value = 0.01960674250039005
$ws.Range("J7").Value = 0.01960674250039005
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.4500266666666666
$ws.Range("N7").Value = 1.35008
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = 1
$ws.Range("Q7").Value = 0.1325649552355556
$ws.Range("R7").Value = 1.19308459712
$ws.Range("S7").Value = 0.01960674250039005
$ws.Range("T7").Value = 0.01960674250039005
